# "one strike in frist frame"
# A new frame (one strike, followed by the usual per-ball counting pattern) is
# inserted at rows 12-14 of the bowling scorecard, pushing the existing frames
# (previously at rows 12-25) down by three rows (to rows 15-28).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new blank rows, shifting rows 12:25 down to 15:28.
$ws.Rows("12:14").Insert()

# The newly inserted rows come in with a generic "new row" style. Re-apply the
# correct look by copying formats from rows that already carry the desired
# style: a "throws" row + "running total" row (both green-filled, s=7) for the
# new rows 12/13, and a blank separator row (s=5) for new row 14.
$ws.Range("A21:W22").Copy()
$ws.Range("A12:W13").PasteSpecial(-4122)
$ws.Range("A20:W20").Copy()
$ws.Range("A14:W14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 12: first ball is a strike ("X"); remaining balls cycle 1-5 like the
# other frames on the sheet.
$ws.Range("A12").Value = "X"
$ws.Range("C12").Value = 3
$ws.Range("D12").Value = 4
$ws.Range("E12").Value = 5
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 2
$ws.Range("H12").Value = 3
$ws.Range("I12").Value = 4
$ws.Range("J12").Value = 5
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 2
$ws.Range("M12").Value = 3
$ws.Range("N12").Value = 4
$ws.Range("O12").Value = 5
$ws.Range("P12").Value = 1
$ws.Range("Q12").Value = 2
$ws.Range("R12").Value = 3
$ws.Range("S12").Value = 4
$ws.Range("T12").Value = 5

# Row 13: running totals for the new frame.
$ws.Range("B13").Value = 17
$ws.Range("C13").Value = 3
$ws.Range("D13").Value = 4
$ws.Range("E13").Value = 5
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 2
$ws.Range("H13").Value = 3
$ws.Range("I13").Value = 4
$ws.Range("J13").Value = 5
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 2
$ws.Range("M13").Value = 3
$ws.Range("N13").Value = 4
$ws.Range("O13").Value = 5
$ws.Range("P13").Value = 1
$ws.Range("Q13").Value = 2
$ws.Range("R13").Value = 3
$ws.Range("S13").Value = 4
$ws.Range("T13").Value = 5
$ws.Range("W13").Value = 74

# Match the author's final cursor position.
$ws.Range("X16").Select()
